$wb = $excel.ActiveWorkbook

# --- Sheet: Trends Status ---
$ws1 = $wb.Worksheets.Item("Trends Status")
$ws1.Range("B2").Value = 1
$ws1.Range("C2").Value = 20
$ws1.Range("D2").Value = 2.6
$ws1.Range("E2").Value = 19.2

$ws1.Range("B3").Value = 3
$ws1.Range("C3").Value = 32
$ws1.Range("D3").Value = 7.9
$ws1.Range("E3").Value = 30.8

$ws1.Range("B4").Value = 18
$ws1.Range("C4").Value = 36
$ws1.Range("D4").Value = 47.4
$ws1.Range("E4").Value = 34.6

$ws1.Range("C5").Value = 8
$ws1.Range("D5").Value = 13.2
$ws1.Range("E5").Value = 7.7

$ws1.Range("B6").Value = 11
$ws1.Range("C6").Value = 8
$ws1.Range("D6").Value = 28.9
$ws1.Range("E6").Value = 7.7

$ws1.Range("B7").Value = 73
$ws1.Range("C7").Value = 150

$ws1.Range("B8").Value = 349
$ws1.Range("C8").Value = 206

# --- Sheet: Priority Status ---
$ws3 = $wb.Worksheets.Item("Priority Status")
$ws3.Range("B2").Value = 103
$ws3.Range("B3").Value = 286
$ws3.Range("B4").Value = 554

# --- Sheet: Species qualification ---
$ws4 = $wb.Worksheets.Item("Species qualification")
$ws4.Range("A2").Value = "SoIB Assessment"
$ws4.Range("B2").Value = 460

$ws4.Range("B3").Value = 111
$ws4.Range("C3").Value = 38

$ws4.Range("C4").Value = 104

# --- Sheet: High Priority break-up -> rename, update values, and duplicate ---
$ws5 = $wb.Worksheets.Item("High Priority break-up")

# Before changing the values, create the new "Major update" sheet as a verbatim
# copy of the current (pre-edit) data, placed right after the renamed sheet.
$ws5.Copy($null, $ws5)
$ws6 = $wb.Worksheets.Item($ws5.Index + 1)
$ws6.Name = "Major update - High Priority "

# Now rename the original sheet and update its values to the "interannual" figures.
$ws5.Name = "Interannual update - High Pri"

$ws5.Range("B2").Value = 73
$ws5.Range("C2").Value = 70.90000000000001
$ws5.Range("D2").Value = 73
$ws5.Range("E2").Value = 84.90000000000001

$ws5.Range("B3").Value = 30
$ws5.Range("C3").Value = 29.1
$ws5.Range("D3").Value = 13
$ws5.Range("E3").Value = 15.1
